$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove old rows 89-123 (the Level 6 / Level 7 event log section being rewritten)
$ws.Range("89:123").Delete()

# Insert 44 fresh blank rows (89-132) to host the updated event log
$ws.Range("89:132").Insert()

$ws.Cells.Item(89,1).Value = 'Y = 800'

$ws.Cells.Item(91,1).Value = 'Enter ship'
$ws.Cells.Item(91,2).Value = 81636
$ws.Cells.Item(91,4).Formula = '=IF(B91="","-",IF(C91="","-",B91-C91))'
$ws.Cells.Item(91,5).Value = '8th'
$ws.Cells.Item(91,6).Value = '16th'
$ws.Cells.Item(91,7).Value = 'triplet'
$ws.Cells.Item(91,8).Value = '8th'
$ws.Rows.Item(91).RowHeight = 16.5

$ws.Rows.Item(92).RowHeight = 16.5

$ws.Cells.Item(93,1).Value = 'end elevator'
$ws.Cells.Item(93,2).Value = 82955
$ws.Cells.Item(93,4).Formula = '=IF(B93="","-",IF(C93="","-",B93-C93))'

$ws.Cells.Item(94,1).Value = 'horse appears'
$ws.Cells.Item(94,2).Value = 83730
$ws.Cells.Item(94,3).Value = 77216
$ws.Cells.Item(94,4).Formula = '=IF(B94="","-",IF(C94="","-",B94-C94))'

$ws.Cells.Item(95,1).Value = 'horse appears'
$ws.Cells.Item(95,2).Value = 86049

$ws.Cells.Item(96,1).Value = 'Level 6'
$rng = $ws.Range("A96:D96")
$rng.Font.Bold = $true
$rng.Font.Italic = $true
$rng.Font.Size = 16
$rng.HorizontalAlignment = -4108
$ws.Rows.Item(96).RowHeight = 20.25

$ws.Cells.Item(97,1).Value = 'horse appears'
$ws.Cells.Item(97,2).Value = 86049
$ws.Cells.Item(97,3).Value = 79536
$ws.Cells.Item(97,4).Formula = '=IF(B97="","-",IF(C97="","-",B97-C97))'

$ws.Cells.Item(98,1).Value = 'q1 appears'
$ws.Cells.Item(98,2).Value = 87616
$ws.Cells.Item(98,4).Formula = '=IF(B98="","-",IF(C98="","-",B98-C98))'

$ws.Cells.Item(99,1).Value = 'q2 appears'
$ws.Cells.Item(99,2).Value = 89186
$ws.Cells.Item(99,3).Value = 82673
$ws.Cells.Item(99,4).Formula = '=IF(B99="","-",IF(C99="","-",B99-C99))'

$ws.Cells.Item(100,1).Value = 'q3 appears'
$ws.Cells.Item(100,2).Value = 92818
$ws.Cells.Item(100,3).Value = 86305
$ws.Cells.Item(100,4).Formula = '=IF(B100="","-",IF(C100="","-",B100-C100))'

$ws.Cells.Item(101,1).Value = 'q4 appears'
$ws.Cells.Item(101,2).Value = 94242
$ws.Cells.Item(101,3).Value = 87729
$ws.Cells.Item(101,4).Formula = '=IF(B101="","-",IF(C101="","-",B101-C101))'

$ws.Cells.Item(102,1).Value = 'q5 appears'
$ws.Cells.Item(102,2).Value = 95650
$ws.Cells.Item(102,3).Value = 89137
$ws.Cells.Item(102,4).Formula = '=IF(B102="","-",IF(C102="","-",B102-C102))'

$ws.Cells.Item(103,1).Value = 'begin exit'
$ws.Cells.Item(103,2).Value = 101844
$ws.Cells.Item(103,3).Value = 95331
$ws.Cells.Item(103,4).Formula = '=IF(B103="","-",IF(C103="","-",B103-C103))'

$ws.Cells.Item(104,1).Value = 'enter station'
$ws.Cells.Item(104,2).Value = 102561
$ws.Cells.Item(104,3).Value = 96048
$ws.Cells.Item(104,4).Formula = '=IF(B104="","-",IF(C104="","-",B104-C104))'

$ws.Cells.Item(105,1).Value = 'black screen'
$ws.Cells.Item(105,2).Value = 103048
$ws.Cells.Item(105,3).Value = 96535
$ws.Cells.Item(105,4).Formula = '=IF(B105="","-",IF(C105="","-",B105-C105))'

$ws.Cells.Item(106,1).Value = 'bonus screen end'
$ws.Cells.Item(106,2).Value = 105199
$ws.Cells.Item(106,3).Value = 98049
$ws.Cells.Item(106,4).Formula = '=IF(B106="","-",IF(C106="","-",B106-C106))'

$ws.Cells.Item(107,1).Value = 'level appears'
$ws.Cells.Item(107,2).Value = 106508
$ws.Cells.Item(107,3).Value = 99353
$ws.Cells.Item(107,4).Formula = '=IF(B107="","-",IF(C107="","-",B107-C107))'

$ws.Cells.Item(108,1).Value = 'Level 7'
$rng = $ws.Range("A108:D108")
$rng.Font.Bold = $true
$rng.Font.Italic = $true
$rng.Font.Size = 16
$rng.HorizontalAlignment = -4108
$ws.Rows.Item(108).RowHeight = 20.25

$ws.Cells.Item(109,1).Value = 'level appears'
$ws.Cells.Item(109,2).Value = 106508
$ws.Cells.Item(109,3).Value = 99353
$ws.Cells.Item(109,4).Formula = '=IF(B109="","-",IF(C109="","-",B109-C109))'

$ws.Cells.Item(110,1).Value = 'enter door'
$ws.Cells.Item(110,2).Value = 107587
$ws.Cells.Item(110,4).Formula = '=IF(B110="","-",IF(C110="","-",B110-C110))'

$ws.Cells.Item(111,1).Value = 'enter door'
$ws.Cells.Item(111,2).Value = 108015
$ws.Cells.Item(111,4).Formula = '=IF(B111="","-",IF(C111="","-",B111-C111))'

$ws.Cells.Item(112,1).Value = 'X = 811'
$ws.Cells.Item(112,2).Value = 108724
$ws.Cells.Item(112,4).Formula = '=IF(B112="","-",IF(C112="","-",B112-C112))'

$ws.Cells.Item(113,1).Value = 'X = 1101'
$ws.Cells.Item(113,2).Value = 108885

$ws.Cells.Item(114,1).Value = 'enter door'
$ws.Cells.Item(114,2).Value = 109301
$ws.Cells.Item(114,4).Formula = '=IF(B114="","-",IF(C114="","-",B114-C114))'

$ws.Cells.Item(115,1).Value = 'enter door'
$ws.Cells.Item(115,2).Value = 109728
$ws.Cells.Item(115,4).Formula = '=IF(B115="","-",IF(C115="","-",B115-C115))'

$ws.Cells.Item(116,1).Value = 'enter door'
$ws.Cells.Item(116,2).Value = 111196
$ws.Cells.Item(116,4).Formula = '=IF(B116="","-",IF(C116="","-",B116-C116))'

$ws.Cells.Item(117,1).Value = 'enter elevator'
$ws.Cells.Item(117,2).Value = 111524
$ws.Cells.Item(117,4).Formula = '=IF(B117="","-",IF(C117="","-",B117-C117))'

$ws.Cells.Item(118,1).Value = 'enter door'
$ws.Cells.Item(118,2).Value = 112599
$ws.Cells.Item(118,4).Formula = '=IF(B118="","-",IF(C118="","-",B118-C118))'

$ws.Cells.Item(119,1).Value = 'enter door'
$ws.Cells.Item(119,2).Value = 113555
$ws.Cells.Item(119,4).Formula = '=IF(B119="","-",IF(C119="","-",B119-C119))'

$ws.Cells.Item(120,1).Value = 'enter door as tempo'
$ws.Cells.Item(120,2).Value = 117635
$ws.Cells.Item(120,4).Formula = '=IF(B120="","-",IF(C120="","-",B120-C120))'

$ws.Cells.Item(121,1).Value = '…'

$ws.Cells.Item(122,1).Value = 'enter elevator'
$ws.Cells.Item(122,2).Value = 121599
$ws.Cells.Item(122,4).Formula = '=IF(B122="","-",IF(C122="","-",B122-C122))'

$ws.Cells.Item(123,1).Value = 'enter top'
$ws.Cells.Item(123,2).Value = 126226
$ws.Cells.Item(123,4).Formula = '=IF(B123="","-",IF(C123="","-",B123-C123))'

$ws.Cells.Item(124,1).Value = 'end chicken fight'
$ws.Cells.Item(124,2).Value = 136739
$ws.Cells.Item(124,4).Formula = '=IF(B124="","-",IF(C124="","-",B124-C124))'

$ws.Cells.Item(125,1).Value = 'screen end'
$ws.Cells.Item(125,2).Value = 140173
$ws.Cells.Item(125,4).Formula = '=IF(B125="","-",IF(C125="","-",B125-C125))'

$ws.Cells.Item(126,1).Value = 'leave windmill room'
$ws.Cells.Item(126,2).Value = 141695
$ws.Cells.Item(126,4).Formula = '=IF(B126="","-",IF(C126="","-",B126-C126))'

$ws.Cells.Item(127,1).Value = 'end crab fight'
$ws.Cells.Item(127,2).Value = 146518
$ws.Cells.Item(127,4).Formula = '=IF(B127="","-",IF(C127="","-",B127-C127))'

$ws.Cells.Item(128,1).Value = 'end riding hood fight'
$ws.Cells.Item(128,2).Value = 150651
$ws.Cells.Item(128,4).Formula = '=IF(B128="","-",IF(C128="","-",B128-C128))'

$ws.Cells.Item(129,1).Value = 'end unfortunate door fight'
$ws.Cells.Item(129,2).Value = 156600
$ws.Cells.Item(129,4).Formula = '=IF(B129="","-",IF(C129="","-",B129-C129))'

$ws.Cells.Item(130,1).Value = 'go up'
$ws.Cells.Item(130,2).Value = 160147
$ws.Cells.Item(130,4).Formula = '=IF(B130="","-",IF(C130="","-",B130-C130))'

$ws.Cells.Item(131,1).Value = 'end level'
$ws.Cells.Item(131,2).Value = 175028
$ws.Cells.Item(131,4).Formula = '=IF(B131="","-",IF(C131="","-",B131-C131))'

$ws.Cells.Item(132,1).Value = 'boss fight end (white screen)'
$ws.Cells.Item(132,2).Value = 179257
$ws.Cells.Item(132,4).Formula = '=IF(B132="","-",IF(C132="","-",B132-C132))'

# Update sheet view: activate, position selection to match new authored state
$ws.Activate()
$ws.Range("C112").Select()
